# Ajustes en la descripción de la evaluación
#
# - Replace the evaluation activity description in row 30 (column J) with the
#   new, shorter wording used by the author.
# - Correct the numbering typo in column H (rows 12-30 had an extra trailing
#   zero: 100,110,...,280 instead of the intended 10,11,...,28 continuing the
#   sequence started in rows 3-11).
# - Leave the selection on the cell that was edited (J30).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ESCALETA")

# 1) Update the evaluation description text for row 30.
$ws.Range("J30").Value = "Evalúa tus conocimientos acerca del tema Los fundamentos de genética"

# 2) Fix the column H values for rows 12-30 (remove the stray trailing zero).
for ($r = 12; $r -le 30; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value = $current / 10
    }
}

# 3) Update the active selection to reflect the cell that was edited.
$ws.Activate()
$ws.Range("J30").Select()
